$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 so everything from the old row 8 lunch/afternoon
# block shifts down by one, matching the "6 hours by turn" schedule fix.
$ws.Rows("9:9").Insert()

# --- Column A: time labels (rebuild the full time column 2..17) ---
$times = @("7:00","7:50","8:40","9:30","9:50","10:40","11:30","12:20","13:00","13:50","14:40","15:30","15:50","16:40","17:30","18:20")
for ($i = 0; $i -lt $times.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $times[$i]
}

# --- Row 3: class moved out of segunda (B3) ---
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

# --- Row 4: class moved from quarta (D4) to quinta (E4) ---
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "MCT-3A-Motores de aplicação"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "ELT-3A-Motores de aplicação"
$ws.Range("F4").Value = "-"

# --- Row 6: class added on quinta (E6) ---
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "MCT-3A-Motores de aplicação"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "ELT-3A-Motores de aplicação"
$ws.Range("F6").Value = "-"

# --- Row 7 unchanged ---
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"

# --- Row 8: now a regular "-" row (lunch moved to row 9) ---
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

# --- Row 9: lunch ("Almoço") block, now at 12:20 ---
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

# --- Rows 10-12: "-" rows ---
foreach ($r in 10..12) {
    $ws.Cells.Item($r, 2).Value = "-"
    $ws.Cells.Item($r, 3).Value = "-"
    $ws.Cells.Item($r, 4).Value = "-"
    $ws.Cells.Item($r, 5).Value = "-"
    $ws.Cells.Item($r, 6).Value = "-"
}

# --- Row 13: "Intervalo" block ---
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

# --- Rows 14-16: "-" rows ---
foreach ($r in 14..16) {
    $ws.Cells.Item($r, 2).Value = "-"
    $ws.Cells.Item($r, 3).Value = "-"
    $ws.Cells.Item($r, 4).Value = "-"
    $ws.Cells.Item($r, 5).Value = "-"
    $ws.Cells.Item($r, 6).Value = "-"
}

# --- Row 17: new trailing row, time only, rest blank (but present) cells ---
foreach ($col in "B17", "C17", "D17", "E17", "F17") {
    $ws.Range($col).Value = ""
    $ws.Range($col).Style = "Normal"
}
